# List of Tasks - Notes.docx
#
# Semantic changes applied (everything else in the document, including
# paragraph order/count/numbering, is left untouched):
#
#   1. "Abstract - Kai (1st Draft Completed)"
#        -> "Abstract - Kai (1st Draft Completed, need to prep for Godshalk)"
#   2. "Hardware in Context - Kai (WIP)"
#        -> "Hardware in Context - Kai (1st Draft Completed, need to prep
#            for Godshalk)"   [with "st" superscripted, matching Abstract]
#   3. "Hardware Description - Kai "
#        -> "Hardware Description - Kai (WIP)"
#
# A few other paragraphs ("Research other HardwareX Papers",
# "Add the HardwareX LaTeX Files", "Setup Zenodo with eDNA GitHub ...")
# simply had their mid-word spell-check markers (w:proofErr) cleared and
# runs re-merged with no visible text change; a self Find/Replace on the
# full phrase reproduces that the same way Word itself would.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Proofing-only touch-ups (text unchanged; clears w:proofErr splits)
# ---------------------------------------------------------------------
$touch = $d.Content
$touch.Find.Execute("Research other HardwareX Papers", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Research other HardwareX Papers", 2) | Out-Null

$touch = $d.Content
$touch.Find.Execute("Add the HardwareX LaTeX Files", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Add the HardwareX LaTeX Files", 2) | Out-Null

$touch = $d.Content
$touch.Find.Execute("Setup Zenodo with eDNA GitHub (Maybe on the UI and Framework GitHub as well)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Setup Zenodo with eDNA GitHub (Maybe on the UI and Framework GitHub as well)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Abstract paragraph: append ", need to prep for Godshalk" before ")"
# ---------------------------------------------------------------------
$abs = $d.Content
$abs.Find.Execute("Draft Completed)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Draft Completed", 2) | Out-Null
$abs.Collapse(0)
$abs.InsertAfter(", need to prep for Godshalk")
$abs.Collapse(0)
$abs.InsertAfter(")")

# ---------------------------------------------------------------------
# 3. Hardware in Context paragraph: replace "(WIP)" with the same
#    "(1st Draft Completed, need to prep for Godshalk)" phrase, with the
#    "st" superscripted like it is in the Abstract line above.
# ---------------------------------------------------------------------
$hic = $d.Content
$hic.Find.Execute("(WIP)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hic.Text = ""

$hic.Collapse(0)
$hic.InsertAfter("(")

$hic.Collapse(0)
$hic.InsertAfter("1")

$hic.Collapse(0)
$hic.InsertAfter("st")
$hic.Font.Superscript = $true

$hic.Collapse(0)
$hic.InsertAfter(" Draft Completed")

$hic.Collapse(0)
$hic.InsertAfter(", ")

$hic.Collapse(0)
$hic.InsertAfter("need to prep for Godshalk")

$hic.Collapse(0)
$hic.InsertAfter(")")

# ---------------------------------------------------------------------
# 4. Hardware Description paragraph: the "(WIP)" marker moved down to
#    here from the "Hardware in Context" line above.
# ---------------------------------------------------------------------
$hdPara = $d.Paragraphs.Item(24)
$hd = $d.Range($hdPara.Range.Start, $hdPara.Range.End - 1)
$hd.Collapse(0)
$hd.InsertAfter("(WIP)")
